$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customers")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "quynh"
$ws.Cells.Item($newRow, 2).Value = "dinh"
$ws.Cells.Item($newRow, 3).Value = "quynh@domain.com"
